$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")
$ws.Range("A2").Copy()
$ws.Range("Z61").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("Z62").PasteSpecial(-4122)
$ws.Range("B6").Copy()
$ws.Range("Z63").PasteSpecial(-4122)
